$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear everything first so the shared-strings table rebuilds fresh, in the exact
# assignment order we use below (mirrors the order produced by the authors' export script).
$ws.Range("A1:T10").ClearContents()

# --- Header row (row 1) : establishes shared-string indices 0-19 ---
$ws.Cells.Item(1, 1).Value() = "Sending cluster"
$ws.Cells.Item(1, 2).Value() = "Ligand symbol"
$ws.Cells.Item(1, 3).Value() = "Receptor symbol"
$ws.Cells.Item(1, 4).Value() = "Target cluster"
$ws.Cells.Item(1, 5).Value() = "Ligand-expressing cells"
$ws.Cells.Item(1, 6).Value() = "Ligand detection rate"
$ws.Cells.Item(1, 7).Value() = "Ligand average expression value"
$ws.Cells.Item(1, 8).Value() = "Ligand total expression value"
$ws.Cells.Item(1, 9).Value() = "Ligand derived specificity of average expression value"
$ws.Cells.Item(1, 10).Value() = "Ligand derived specificity of total expression value"
$ws.Cells.Item(1, 11).Value() = "Receptor-expressing cells"
$ws.Cells.Item(1, 12).Value() = "Receptor detection rate"
$ws.Cells.Item(1, 13).Value() = "Receptor average expression value"
$ws.Cells.Item(1, 14).Value() = "Receptor total expression value"
$ws.Cells.Item(1, 15).Value() = "Receptor derived specificity of average expression value"
$ws.Cells.Item(1, 16).Value() = "Receptor derived specificity of total expression value"
$ws.Cells.Item(1, 17).Value() = "Edge average expression weight"
$ws.Cells.Item(1, 18).Value() = "Edge total expression weight"
$ws.Cells.Item(1, 19).Value() = "Edge average expression derived specificity"
$ws.Cells.Item(1, 20).Value() = "Edge total expression derived specificity"

# --- String columns for data rows 2-7, written column-by-column so the shared-strings
#     table is rebuilt in the order: FAPs, MuSCs, Efna5, Epha1, ECs (indices 20-24) ---
# Column A
$ws.Cells.Item(2, 1).Value() = "FAPs"
$ws.Cells.Item(3, 1).Value() = "FAPs"
$ws.Cells.Item(4, 1).Value() = "FAPs"
$ws.Cells.Item(5, 1).Value() = "MuSCs"
$ws.Cells.Item(6, 1).Value() = "MuSCs"
$ws.Cells.Item(7, 1).Value() = "MuSCs"

# Column B
$ws.Cells.Item(2, 2).Value() = "Efna5"
$ws.Cells.Item(3, 2).Value() = "Efna5"
$ws.Cells.Item(4, 2).Value() = "Efna5"
$ws.Cells.Item(5, 2).Value() = "Efna5"
$ws.Cells.Item(6, 2).Value() = "Efna5"
$ws.Cells.Item(7, 2).Value() = "Efna5"

# Column C
$ws.Cells.Item(2, 3).Value() = "Epha1"
$ws.Cells.Item(3, 3).Value() = "Epha1"
$ws.Cells.Item(4, 3).Value() = "Epha1"
$ws.Cells.Item(5, 3).Value() = "Epha1"
$ws.Cells.Item(6, 3).Value() = "Epha1"
$ws.Cells.Item(7, 3).Value() = "Epha1"

# Column D
$ws.Cells.Item(2, 4).Value() = "ECs"
$ws.Cells.Item(3, 4).Value() = "FAPs"
$ws.Cells.Item(4, 4).Value() = "MuSCs"
$ws.Cells.Item(5, 4).Value() = "ECs"
$ws.Cells.Item(6, 4).Value() = "FAPs"
$ws.Cells.Item(7, 4).Value() = "MuSCs"

# --- Numeric columns for data rows 2-7 ---
# Row 2
$ws.Cells.Item(2, 5).Value() = 3
$ws.Cells.Item(2, 6).Value() = 1
$ws.Cells.Item(2, 7).Value() = 2.900731333333333
$ws.Cells.Item(2, 8).Value() = 8.702194
$ws.Cells.Item(2, 9).Value() = 0.8130494232775288
$ws.Cells.Item(2, 10).Value() = 0.8130494232775289
$ws.Cells.Item(2, 11).Value() = 3
$ws.Cells.Item(2, 12).Value() = 1
$ws.Cells.Item(2, 13).Value() = 1.334383666666667
$ws.Cells.Item(2, 14).Value() = 4.003151
$ws.Cells.Item(2, 15).Value() = 0.1312069045987744
$ws.Cells.Item(2, 16).Value() = 0.1312069045987744
$ws.Cells.Item(2, 17).Value() = 3.870688512588222
$ws.Cells.Item(2, 18).Value() = 34.836196613294
$ws.Cells.Item(2, 19).Value() = 0.1066776981140632
$ws.Cells.Item(2, 20).Value() = 0.1066776981140633

# Row 3
$ws.Cells.Item(3, 5).Value() = 3
$ws.Cells.Item(3, 6).Value() = 1
$ws.Cells.Item(3, 7).Value() = 2.900731333333333
$ws.Cells.Item(3, 8).Value() = 8.702194
$ws.Cells.Item(3, 9).Value() = 0.8130494232775288
$ws.Cells.Item(3, 10).Value() = 0.8130494232775289
$ws.Cells.Item(3, 11).Value() = 3
$ws.Cells.Item(3, 12).Value() = 1
$ws.Cells.Item(3, 13).Value() = 2.477462333333333
$ws.Cells.Item(3, 14).Value() = 7.432386999999999
$ws.Cells.Item(3, 15).Value() = 0.2436032245723858
$ws.Cells.Item(3, 16).Value() = 0.2436032245723858
$ws.Cells.Item(3, 17).Value() = 7.18645261745311
$ws.Cells.Item(3, 18).Value() = 64.678073557078
$ws.Cells.Item(3, 19).Value() = 0.1980614612471246
$ws.Cells.Item(3, 20).Value() = 0.1980614612471247

# Row 4
$ws.Cells.Item(4, 5).Value() = 3
$ws.Cells.Item(4, 6).Value() = 1
$ws.Cells.Item(4, 7).Value() = 2.900731333333333
$ws.Cells.Item(4, 8).Value() = 8.702194
$ws.Cells.Item(4, 9).Value() = 0.8130494232775288
$ws.Cells.Item(4, 10).Value() = 0.8130494232775289
$ws.Cells.Item(4, 11).Value() = 3
$ws.Cells.Item(4, 12).Value() = 1
$ws.Cells.Item(4, 13).Value() = 6.358226000000001
$ws.Cells.Item(4, 14).Value() = 19.074678
$ws.Cells.Item(4, 15).Value() = 0.6251898708288398
$ws.Cells.Item(4, 16).Value() = 0.6251898708288398
$ws.Cells.Item(4, 17).Value() = 18.44350538261467
$ws.Cells.Item(4, 18).Value() = 165.991548443532
$ws.Cells.Item(4, 19).Value() = 0.5083102639163409
$ws.Cells.Item(4, 20).Value() = 0.508310263916341

# Row 5
$ws.Cells.Item(5, 5).Value() = 3
$ws.Cells.Item(5, 6).Value() = 1
$ws.Cells.Item(5, 7).Value() = 0.6669870000000001
$ws.Cells.Item(5, 8).Value() = 2.000961
$ws.Cells.Item(5, 9).Value() = 0.1869505767224711
$ws.Cells.Item(5, 10).Value() = 0.1869505767224711
$ws.Cells.Item(5, 11).Value() = 3
$ws.Cells.Item(5, 12).Value() = 1
$ws.Cells.Item(5, 13).Value() = 1.334383666666667
$ws.Cells.Item(5, 14).Value() = 4.003151
$ws.Cells.Item(5, 15).Value() = 0.1312069045987744
$ws.Cells.Item(5, 16).Value() = 0.1312069045987744
$ws.Cells.Item(5, 17).Value() = 0.8900165586790001
$ws.Cells.Item(5, 18).Value() = 8.010149028111
$ws.Cells.Item(5, 19).Value() = 0.02452920648471111
$ws.Cells.Item(5, 20).Value() = 0.02452920648471111

# Row 6
$ws.Cells.Item(6, 5).Value() = 3
$ws.Cells.Item(6, 6).Value() = 1
$ws.Cells.Item(6, 7).Value() = 0.6669870000000001
$ws.Cells.Item(6, 8).Value() = 2.000961
$ws.Cells.Item(6, 9).Value() = 0.1869505767224711
$ws.Cells.Item(6, 10).Value() = 0.1869505767224711
$ws.Cells.Item(6, 11).Value() = 3
$ws.Cells.Item(6, 12).Value() = 1
$ws.Cells.Item(6, 13).Value() = 2.477462333333333
$ws.Cells.Item(6, 14).Value() = 7.432386999999999
$ws.Cells.Item(6, 15).Value() = 0.2436032245723858
$ws.Cells.Item(6, 16).Value() = 0.2436032245723858
$ws.Cells.Item(6, 17).Value() = 1.652435169323
$ws.Cells.Item(6, 18).Value() = 14.871916523907
$ws.Cells.Item(6, 19).Value() = 0.04554176332526116
$ws.Cells.Item(6, 20).Value() = 0.04554176332526117

# Row 7
$ws.Cells.Item(7, 5).Value() = 3
$ws.Cells.Item(7, 6).Value() = 1
$ws.Cells.Item(7, 7).Value() = 0.6669870000000001
$ws.Cells.Item(7, 8).Value() = 2.000961
$ws.Cells.Item(7, 9).Value() = 0.1869505767224711
$ws.Cells.Item(7, 10).Value() = 0.1869505767224711
$ws.Cells.Item(7, 11).Value() = 3
$ws.Cells.Item(7, 12).Value() = 1
$ws.Cells.Item(7, 13).Value() = 6.358226000000001
$ws.Cells.Item(7, 14).Value() = 19.074678
$ws.Cells.Item(7, 15).Value() = 0.6251898708288398
$ws.Cells.Item(7, 16).Value() = 0.6251898708288398
$ws.Cells.Item(7, 17).Value() = 4.240854085062002
$ws.Cells.Item(7, 18).Value() = 38.16768676555801
$ws.Cells.Item(7, 19).Value() = 0.1168796069124988
$ws.Cells.Item(7, 20).Value() = 0.1168796069124988

# --- Remove the now-empty trailing rows (originally 8-10) so the sheet's dimension
#     shrinks back down to A1:T7 ---
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(8).Delete()
